# Insert a new data row at row 368 (pushing existing rows 368:391 down to
# 369:392) in the "Feria Lagunitas de Puerto Montt - Papa" price listing,
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 368 downward (inserting a new blank row 368).
$ws.Rows("368:368").Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A368").Value = 4
$ws.Range("B368").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C368").Value = "Los Lagos"
$ws.Range("D368").Value = 44706
$ws.Range("E368").Value = 10
$ws.Range("F368").Value = 100114001
$ws.Range("G368").Value = "Papa"
$ws.Range("H368").Value = "Patagonia"
$ws.Range("I368").Value = "1a (guarda)"
$ws.Range("J368").Value = 150
$ws.Range("K368").Value = 7000
$ws.Range("L368").Value = 7500
$ws.Range("M368").Value = 7233
$ws.Range("N368").Value = "`$/saco 25 kilos"
$ws.Range("O368").Value = "Provincia de Llanquihue"
$ws.Range("P368").Value = 289
$ws.Range("Q368").Value = 25
$ws.Range("R368").Value = "Hortaliza"
